$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Set values in the exact order that introduces new shared strings,
# --- so the resulting shared-strings table indices line up with the target.

# 1) D2: "Concluída" -> "Melhorar"
$ws.Range("D2").Value = "Melhorar"

# 2) A7: "gerente" -> "Motorista"
$ws.Range("A7").Value = "Motorista"

# 3) C10: (empty) -> "Cadastra as despesas geradas pelos eventos"
$ws.Range("C10").Value = "Cadastra as despesas geradas pelos eventos"

# 4) New row 14, column B: "Banco"
$ws.Range("B14").Value = "Banco"

# 5) New column E header (E1): "Prioridades"
$ws.Range("E1").Value = "Prioridades"

# 6) New row 14, column A: "Editar / excluir"
$ws.Range("A14").Value = "Editar / excluir"

# 7) New row 14, column C: "Permitir alterações de dados"
$ws.Range("C14").Value = "Permitir alterações de dados"

# --- Remaining value updates (reuse already-created shared strings) ---

# D12: "Concluída" -> "Melhorar"
$ws.Range("D12").Value = "Melhorar"

# New row 14, column D: "Desenvolvimento"
$ws.Range("D14").Value = "Desenvolvimento"

# New column E (Prioridades) numeric values
$ws.Range("E2").Value = 1
$ws.Range("E10").Value = 1
$ws.Range("E11").Value = 1
$ws.Range("E12").Value = 2
$ws.Range("E13").Value = 1
$ws.Range("E14").Value = 2

# --- Formatting: reuse the existing wrap-text style (style index 1) for the
# --- new text cells, matching the style used throughout columns A-D.
$ws.Range("E1").WrapText = $true
$ws.Range("A14").WrapText = $true
$ws.Range("B14").WrapText = $true
$ws.Range("C14").WrapText = $true
$ws.Range("D14").WrapText = $true

# Row 1 becomes taller to fit the new header
$ws.Rows.Item(1).RowHeight = 30

# Final selection, mirroring the saved workbook state
$ws.Range("A15").Select() | Out-Null
